$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells D1:F1, with the same header style as the existing A1:C1 ---
$ws.Range("D1").Value = "USD"
$ws.Range("E1").Value = "EUR"
$ws.Range("F1").Value = "CNY"
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Full data block (rows 2-14, columns A:F) -----------------------------
# Column A holds dates stored as literal text (not date serials), so each
# date string is entered with a leading apostrophe to force text entry.
$data = @(
    @{ Row = 2;  A = "'2023-09-01"; B = 250;  C = "expenses"; D = 96.33;             E = 104.94;            F = 13.19 },
    @{ Row = 3;  A = "'2023-09-05"; B = 500;  C = "expenses"; D = 96.62;             E = 104.42;            F = 13.29 },
    @{ Row = 4;  A = "'2023-09-25"; B = 5000; C = "expenses"; D = 96.04000000000001; E = 102.25;            F = 13.14 },
    @{ Row = 5;  A = "'2023-10-10"; B = 2500; C = "expenses"; D = 101.36;            E = 107.03;            F = 13.89 },
    @{ Row = 6;  A = "'2023-10-15"; B = 565;  C = "expenses"; D = 97.31;             E = 102.55;            F = 13.3  },
    @{ Row = 7;  A = "'2023-10-30"; B = 1000; C = "expenses"; D = 93.22;             E = 98.34999999999999; F = 12.71 },
    @{ Row = 8;  A = "'2023-11-03"; B = 25;   C = "expenses"; D = 93.17;             E = 99;                F = 12.7  },
    @{ Row = 9;  A = "'2023-09-01"; B = 1500; C = "income";   D = 96.33;             E = 104.94;            F = 13.19 },
    @{ Row = 10; A = "'2023-09-02"; B = 500;  C = "income";   D = 96.34;             E = 104.61;            F = 13.25 },
    @{ Row = 11; A = "'2023-09-25"; B = 2500; C = "income";   D = 96.04000000000001; E = 102.25;            F = 13.14 },
    @{ Row = 12; A = "'2023-09-30"; B = 5500; C = "income";   D = 97.41;             E = 103.16;            F = 13.36 },
    @{ Row = 13; A = "'2023-10-01"; B = 1000; C = "income";   D = 97.41;             E = 103.16;            F = 13.36 },
    @{ Row = 14; A = "'2023-10-15"; B = 250;  C = "income";   D = 97.31;             E = 102.55;            F = 13.3  }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
